$d = $word.ActiveDocument
$bullet = [char]0x2022

# ---------------------------------------------------------------------------
# 1. Collapse the three long "CORE COMPETENCIES" paragraphs into a single
#    short summary paragraph.
# ---------------------------------------------------------------------------
$p1 = $null
$p2 = $null
$p3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Data Visualization & Design: Interactive Dashboards*") {
        $p1 = $p
    } elseif ($p.Range.Text -like "Geospatial Analysis & Mapping: Spatial Analysis*") {
        $p2 = $p
    } elseif ($p.Range.Text -like "Technical Visualization: Programming*") {
        $p3 = $p
    }
}

$newSummary = "Data Visualization & Design " + $bullet + " Geospatial Analysis & Mapping " + $bullet + " Technical Visualization"
$p1.Range.Text = $newSummary

$deleteRange = $d.Range($p2.Range.Start, $p3.Range.End)
$deleteRange.Delete()

# ---------------------------------------------------------------------------
# 2. Append a new "TECHNICAL SKILLS" section at the end of the document,
#    before the final section properties.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Range.Text = "TECHNICAL SKILLS"
$headingPara.Style = "Heading 2"

$headingPara.Range.InsertParagraphAfter()
$dataVizPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$dataVizPara.Style = "Normal"
$dataVizPara.Range.Text = "DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design; Web Visualization; Presentation Design; Data Storytelling"

$dataVizPara.Range.InsertParagraphAfter()
$geoPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$geoPara.Style = "Normal"
$geoPara.Range.Text = "GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing; Census Data Integration; Custom Tile Servers; Spatial Clustering"

$geoPara.Range.InsertParagraphAfter()
$techPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$techPara.Style = "Normal"
$techPara.Range.Text = "TECHNICAL VISUALIZATION Programming; Database Integration; Cloud Platforms; Web Technologies; Statistical Computing; Version Control; DevOps"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
